# Train Run Trends - append the new day's row (2016-05-24, serial 42514)
# to the "Data" sheet, following the same column layout as every prior row:
#   A: Date
#   B: Total PTC Runs
#   C: Single Init Runs
#   D: Multiple Init Runs
#   E: Cut Out Runs
#   F: Total Completed
#   G: Total Completed %   (plain value, matching rows 37-40 which no longer carry a formula)
#   H: Completed Trip Length Average
#   I: Completed Trip Length Min
#   J: Completed Trip Length Max

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("A41").Value = 42514
$ws.Range("B41").Value = 117
$ws.Range("C41").Value = 91
$ws.Range("D41").Value = 1
$ws.Range("E41").Value = 25
$ws.Range("F41").Value = 92
$ws.Range("G41").Value = 0.78632478632478631
$ws.Range("H41").Value = 42.115099714529642
$ws.Range("I41").Value = 25.833333337213844
$ws.Range("J41").Value = 193.56666667386889

# Leave the new row selected (matches the author's last on-sheet action).
$ws.Range("H41:J41").Select() | Out-Null
